$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.721.69'
$ws.Cells.Item(2, 5).Value = '  +0.09%  '
$ws.Cells.Item(3, 4).Value = '1.902.39'
$ws.Cells.Item(3, 5).Value = '  +0.37%  '
$ws.Cells.Item(4, 4).Value = "'0.9990"
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  -0.19%  '
$ws.Cells.Item(5, 4).Value = "'312.02"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -0.18%  '
$ws.Cells.Item(6, 4).Value = "'0.9978"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -0.29%  '
$ws.Cells.Item(7, 4).Value = "'0.5218"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  +6.54%  '
$ws.Cells.Item(8, 4).Value = "'0.3776"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -0.50%  '
$ws.Cells.Item(9, 4).Value = "'0.07221"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -1.52%  '
$ws.Cells.Item(10, 4).Value = "'21.27"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +3.39%  '
$ws.Cells.Item(11, 4).Value = "'0.9076"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -0.73%  '
$ws.Cells.Item(12, 2).Value = 'WrappedEther'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(12, 4).Value = '1.925.69'
$ws.Cells.Item(12, 5).Value = '  -0.35%  '
$ws.Cells.Item(13, 2).Value = 'TRON'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(13, 4).Value = "'0.07621"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -0.82%  '
$ws.Cells.Item(14, 4).Value = "'5.440"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -0.69%  '
$ws.Cells.Item(15, 4).Value = "'91.97"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +0.99%  '
$ws.Cells.Item(16, 4).Value = "'0.9977"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -0.35%  '
$ws.Cells.Item(17, 4).Value = "'0.000008672"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -1.23%  '
$ws.Cells.Item(18, 4).Value = "'0.9977"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -0.30%  '
$ws.Cells.Item(19, 4).Value = '27.769.33'
$ws.Cells.Item(19, 5).Value = '  +0.14%  '
$ws.Cells.Item(20, 4).Value = "'14.50"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -0.01%  '
$ws.Cells.Item(21, 4).Value = "'5.140"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +0.27%  '
$ws.Cells.Item(22, 4).Value = '2.155.76'
$ws.Cells.Item(22, 5).Value = '  +0.85%  '
$ws.Cells.Item(23, 4).Value = "'10.83"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +0.74%  '
$ws.Cells.Item(24, 4).Value = "'6.594"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -0.26%  '
$ws.Cells.Item(25, 4).Value = "'153.20"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -0.31%  '
$ws.Cells.Item(26, 4).Value = "'1.868"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -1.88%  '
$ws.Cells.Item(27, 4).Value = "'2.161"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +0.08%  '
$ws.Cells.Item(28, 4).Value = "'18.28"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -0.51%  '
$ws.Cells.Item(29, 4).Value = "'114.29"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -1.23%  '
$ws.Cells.Item(30, 4).Value = "'4.835"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -1.45%  '
$ws.Cells.Item(31, 4).Value = "'0.08983"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +0.80%  '
$ws.Cells.Item(32, 4).Value = "'4.876"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +5.04%  '
$ws.Cells.Item(33, 4).Value = "'3.172"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -0.65%  '
$ws.Cells.Item(34, 4).Value = "'1.226"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +0.14%  '
$ws.Cells.Item(35, 4).Value = "'0.7753"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +1.30%  '
$ws.Cells.Item(36, 4).Value = "'2.622"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +3.40%  '
$ws.Cells.Item(37, 5).Value = '  +2.55%  '
$ws.Cells.Item(38, 4).Value = "'3.068"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +2.92%  '
$ws.Cells.Item(39, 4).Value = "'1.091"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -0.46%  '
$ws.Cells.Item(40, 4).Value = "'0.5515"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +0.54%  '
$ws.Cells.Item(41, 4).Value = "'0.05266"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -0.27%  '
$ws.Cells.Item(42, 4).Value = "'6.658"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -3.75%  '
$ws.Cells.Item(43, 4).Value = "'114.29"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +2.93%  '
$ws.Cells.Item(44, 4).Value = "'8.507"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -0.18%  '
$ws.Cells.Item(45, 4).Value = "'0.1509"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -0.59%  '
$ws.Cells.Item(46, 4).Value = "'0.4795"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -0.12%  '
$ws.Cells.Item(47, 4).Value = "'10.48"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -1.74%  '
$ws.Cells.Item(48, 4).Value = "'0.9969"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -0.37%  '
$ws.Cells.Item(49, 4).Value = "'1.617"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -1.00%  '
$ws.Cells.Item(50, 4).Value = "'66.73"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -1.20%  '
$ws.Cells.Item(51, 4).Value = "'0.05992"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -0.94%  '
